$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 315.85715
$ws.Range("I28").Value = 198.16667
$ws.Range("K28").Value = 198.16667
$ws.Range("M28").Value = 286.83333
$ws.Range("H69").Value = 1534.2106
$ws.Range("J69").Value = 1491.6666
$ws.Range("L69").Value = 4474.9998
$ws.Range("N69").Value = -6222.9998
$ws.Range("H72").Value = 1534.2106
$ws.Range("J72").Value = 1491.6666
$ws.Range("L72").Value = 13424.9994
$ws.Range("N72").Value = -22160.9994
$ws.Range("H76").Value = 2317993.2
$ws.Range("I76").Value = 3293.3333
$ws.Range("K76").Value = 3293.3333
$ws.Range("M76").Value = -2978.3333
$ws.Range("H79").Value = 2317993.2
$ws.Range("I79").Value = 3293.3333
$ws.Range("K79").Value = 3293.3333
$ws.Range("M79").Value = -2201.3333
$ws.Range("H129").Value = 209051.48
$ws.Range("I129").Value = 405.875
$ws.Range("J129").Value = 250780.6
$ws.Range("K129").Value = 1217.625
$ws.Range("L129").Value = 752341.8
$ws.Range("M129").Value = 3782.375
$ws.Range("N129").Value = -762341.8
$ws.Range("H137").Value = 1635.5
$ws.Range("I137").Value = 1686.3158
$ws.Range("K137").Value = 5058.9474
$ws.Range("M137").Value = -2508.9474
$ws.Range("H138").Value = 2182.4849
$ws.Range("I138").Value = 1276.9524
$ws.Range("J138").Value = 3767.1667
$ws.Range("K138").Value = 3830.857199999999
$ws.Range("L138").Value = 11301.5001
$ws.Range("M138").Value = 1309.142800000001
$ws.Range("N138").Value = -21581.5001
$ws.Range("H141").Value = 3025.3635
$ws.Range("I141").Value = 2810.5715
$ws.Range("J141").Value = 3401.25
$ws.Range("K141").Value = 8431.7145
$ws.Range("L141").Value = 10203.75
$ws.Range("M141").Value = -3251.7145
$ws.Range("N141").Value = -20563.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 935.5789
$ws.Range("I2").Value = 1027.3572
$ws.Range("J2").Value = 678.6
$ws.Range("K2").Value = 1027.3572
$ws.Range("L2").Value = 678.6
$ws.Range("M2").Value = -914.3571999999999
$ws.Range("N2").Value = -904.6
$ws.Range("H32").Value = 6745.2812
$ws.Range("I32").Value = 5861.3125
$ws.Range("K32").Value = 5861.3125
$ws.Range("M32").Value = -5574.3125
$ws.Range("H45").Value = 3055.5
$ws.Range("I45").Value = 3179.923
$ws.Range("J45").Value = 2978.476
$ws.Range("K45").Value = 3179.923
$ws.Range("L45").Value = 2978.476
$ws.Range("M45").Value = -2802.923
$ws.Range("N45").Value = -3732.476
$ws.Range("H61").Value = 2962.7058
$ws.Range("I61").Value = 2883.923
$ws.Range("J61").Value = 3218.75
$ws.Range("K61").Value = 2883.923
$ws.Range("L61").Value = 3218.75
$ws.Range("M61").Value = -2671.923
$ws.Range("N61").Value = -3642.75
$ws.Range("H63").Value = 6252171
$ws.Range("I63").Value = 2714
$ws.Range("J63").Value = 31250000
$ws.Range("K63").Value = 2714
$ws.Range("L63").Value = 31250000
$ws.Range("M63").Value = -2028
$ws.Range("N63").Value = -31251372
$ws.Range("H66").Value = 6252171
$ws.Range("I66").Value = 2714
$ws.Range("J66").Value = 31250000
$ws.Range("K66").Value = 13570
$ws.Range("L66").Value = 156250000
$ws.Range("M66").Value = -10138
$ws.Range("N66").Value = -156256864
$ws.Range("H116").Value = 935.5789
$ws.Range("I116").Value = 1027.3572
$ws.Range("J116").Value = 678.6
$ws.Range("K116").Value = 1027.3572
$ws.Range("L116").Value = 678.6
$ws.Range("M116").Value = 1266.6428
$ws.Range("N116").Value = -5266.6
$ws.Range("H136").Value = 2962.7058
$ws.Range("I136").Value = 2883.923
$ws.Range("J136").Value = 3218.75
$ws.Range("K136").Value = 8651.769
$ws.Range("L136").Value = 9656.25
$ws.Range("M136").Value = -6101.769
$ws.Range("N136").Value = -14756.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 935.5789
$ws.Range("I3").Value = 1027.3572
$ws.Range("J3").Value = 678.6
$ws.Range("K3").Value = 1027.3572
$ws.Range("L3").Value = 678.6
$ws.Range("M3").Value = -913.3571999999999
$ws.Range("N3").Value = -906.6
$ws.Range("H20").Value = 2774.8125
$ws.Range("I20").Value = 2907.3333
$ws.Range("K20").Value = 2907.3333
$ws.Range("M20").Value = -2660.3333
$ws.Range("H39").Value = 15000
$ws.Range("J39").Value = 15000
$ws.Range("L39").Value = 15000
$ws.Range("N39").Value = -15778
$ws.Range("H42").Value = 119995
$ws.Range("J42").Value = 119995
$ws.Range("L42").Value = 119995
$ws.Range("N42").Value = -120651
$ws.Range("H105").Value = 1472343
$ws.Range("I105").Value = 1544
$ws.Range("J105").Value = 2633500
$ws.Range("K105").Value = 1544
$ws.Range("L105").Value = 2633500
$ws.Range("M105").Value = 203
$ws.Range("N105").Value = -2636994
$ws.Range("H107").Value = 1342.5
$ws.Range("J107").Value = 1640.8572
$ws.Range("L107").Value = 1640.8572
$ws.Range("N107").Value = -5480.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 18018.6
$ws.Range("I58").Value = 1313.2
$ws.Range("J58").Value = 101545.6
$ws.Range("K58").Value = 1313.2
$ws.Range("L58").Value = 101545.6
$ws.Range("M58").Value = -1110.2
$ws.Range("N58").Value = -101951.6
$ws.Range("H62").Value = 4728.5713
$ws.Range("I62").Value = 4516.6665
$ws.Range("K62").Value = 4516.6665
$ws.Range("M62").Value = -3892.6665
$ws.Range("H65").Value = 4728.5713
$ws.Range("I65").Value = 4516.6665
$ws.Range("K65").Value = 22583.3325
$ws.Range("M65").Value = -19463.3325
$ws.Range("H136").Value = 18018.6
$ws.Range("I136").Value = 1313.2
$ws.Range("J136").Value = 101545.6
$ws.Range("K136").Value = 3939.6
$ws.Range("L136").Value = 304636.8
$ws.Range("M136").Value = -1389.6
$ws.Range("N136").Value = -309736.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 722.98
$ws.Range("J131").Value = 723.9192
$ws.Range("L131").Value = 2171.7576
$ws.Range("N131").Value = -12251.7576
$ws.Range("H132").Value = 860.8
$ws.Range("J132").Value = 995
$ws.Range("L132").Value = 8955
$ws.Range("N132").Value = -14015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 18000
$ws.Range("J15").Value = 18000
$ws.Range("L15").Value = 18000
$ws.Range("N15").Value = -18576
$ws.Range("H39").Value = 29999
$ws.Range("J39").Value = 29999
$ws.Range("L39").Value = 29999
$ws.Range("N39").Value = -31063
$ws.Range("H70").Value = 3133395.2
$ws.Range("I70").Value = 3782.5454
$ws.Range("J70").Value = 6958477.5
$ws.Range("K70").Value = 3782.5454
$ws.Range("L70").Value = 6958477.5
$ws.Range("M70").Value = -3512.5454
$ws.Range("N70").Value = -6959017.5
$ws.Range("H73").Value = 3133395.2
$ws.Range("I73").Value = 3782.5454
$ws.Range("J73").Value = 6958477.5
$ws.Range("K73").Value = 3782.5454
$ws.Range("L73").Value = 6958477.5
$ws.Range("M73").Value = -2846.5454
$ws.Range("N73").Value = -6960349.5
$ws.Range("H80").Value = 2328.2058
$ws.Range("I80").Value = 1275.7222
$ws.Range("J80").Value = 3512.25
$ws.Range("K80").Value = 1275.7222
$ws.Range("L80").Value = 3512.25
$ws.Range("M80").Value = -277.7221999999999
$ws.Range("N80").Value = -5508.25
$ws.Range("H81").Value = 18000
$ws.Range("J81").Value = 18000
$ws.Range("L81").Value = 18000
$ws.Range("N81").Value = -19996
$ws.Range("H83").Value = 2328.2058
$ws.Range("I83").Value = 1275.7222
$ws.Range("J83").Value = 3512.25
$ws.Range("K83").Value = 6378.611
$ws.Range("L83").Value = 17561.25
$ws.Range("M83").Value = -1386.611
$ws.Range("N83").Value = -27545.25
$ws.Range("H84").Value = 18000
$ws.Range("J84").Value = 18000
$ws.Range("L84").Value = 54000
$ws.Range("N84").Value = -63984
$ws.Range("H132").Value = 21396
$ws.Range("I132").Value = 4976.778
$ws.Range("J132").Value = 48263.816
$ws.Range("K132").Value = 14930.334
$ws.Range("L132").Value = 144791.448
$ws.Range("M132").Value = -12400.334
$ws.Range("N132").Value = -149851.448

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1770.8214
$ws.Range("I132").Value = 1118.409
$ws.Range("K132").Value = 3355.227
$ws.Range("M132").Value = -825.2270000000003
$ws.Range("H136").Value = 1268
$ws.Range("I136").Value = 1268
$ws.Range("K136").Value = 3804
$ws.Range("M136").Value = -1254

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1187.5454
$ws.Range("I113").Value = 1278.0526
$ws.Range("J113").Value = 614.3333
$ws.Range("K113").Value = 3834.1578
$ws.Range("L113").Value = 1842.9999
$ws.Range("M113").Value = -1664.1578
$ws.Range("N113").Value = -6182.9999
$ws.Range("H132").Value = 1946.75
$ws.Range("I132").Value = 1652.2
$ws.Range("J132").Value = 2241.3
$ws.Range("K132").Value = 4956.6
$ws.Range("L132").Value = 6723.900000000001
$ws.Range("M132").Value = -2426.6
$ws.Range("N132").Value = -11783.9
